# Build min_em obj function
# Populate the solid_wood (B) and sum_product (D) results for periods 6-10 (rows 7-11)
# and recompute the dependent ecosystem (E) / system (F) columns to match the
# values produced by the (external) minimum-emissions objective function run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (solid_wood, paper, ecosystem)
$updates = @(
    @{ Row = 7;  SolidWood = 6.8425;             Paper = 0; Ecosystem = 149.2866390832741 },
    @{ Row = 8;  SolidWood = 5.416192502625803;  Paper = 0; Ecosystem = 140.2214236083316 },
    @{ Row = 9;  SolidWood = 4.287196379320418;  Paper = 0; Ecosystem = 142.435132491852  },
    @{ Row = 10; SolidWood = 6.130537579387612;  Paper = 0; Ecosystem = 149.2866390832741 },
    @{ Row = 11; SolidWood = 5.2436374387351;    Paper = 0; Ecosystem = 149.2866390832741 }
)

foreach ($u in $updates) {
    $row = $u.Row
    $solidWood = $u.SolidWood
    $paper = $u.Paper
    $ecosystem = $u.Ecosystem
    $sumProduct = $solidWood + $paper
    $system = $ecosystem + $sumProduct

    $ws.Cells.Item($row, 2).Value = $solidWood   # B: solid_wood
    $ws.Cells.Item($row, 3).Value = $paper        # C: paper
    $ws.Cells.Item($row, 4).Value = $sumProduct   # D: sum_product
    $ws.Cells.Item($row, 5).Value = $ecosystem    # E: ecosystem
    $ws.Cells.Item($row, 6).Value = $system       # F: system
}
